# Add the new "2022" mortality-data column (column S) to Supp_Table_2.
#
# The header "2022" looks like a number, and a plain
#   $ws.Range("S1").Value = "2022"
# would be auto-coerced by Excel into the NUMBER 2022 (same as any other
# digit-only string typed into a cell). Every other year header in this
# sheet (2005 ... 2021, column B..R) is stored as TEXT, so we need the new
# header to be text too. Building it with a formula that concatenates two
# string literals forces a text result; copying that result back onto
# itself as values-only collapses the formula down to a plain text literal
# without ever touching the cell's number format (so no stray style gets
# introduced).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S1").Formula = "=""20""&""22"""
$ws.Range("S1").Copy()
$ws.Range("S1").PasteSpecial(-4163)  # xlPasteValues

# 2022 counts for each ICD10 group, rows 2-21 (column S), in sheet order.
$values = @(563, 29, 21, 2188, 1340, 144, 384, 230, 19, 255, 58, 376, 42, 34, 19, 136, 95, 60, 7, 120)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 19).Value = $values[$i]
}
